# The "Förändrad" column (C) holds a date serial that was bumped by one day
# (2023-10-03 -> 2023-10-04, i.e. serial 45202 -> 45203) for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item(1, 3).End(4).Row  # xlDown from header row -> last contiguous data row
if ($lastRow -lt 2) { $lastRow = 494 }

$rng = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3))
foreach ($cell in $rng) {
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
